$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Corrigido localização do Keys:
# the phone number on row 3 (Junior) had the wrong area/location code baked
# into it - correct it to the right number.
$ws.Range("C3").Value = 5532991170287

# The "numero" column header (C1) and the corrected number (C3) should use
# the same centered integer format that the first phone number (C2) already
# uses, so the whole "numero" column reads consistently.
$ws.Range("C1").NumberFormat = "0"
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4108

$ws.Range("C3").NumberFormat = "0"
$ws.Range("C3").HorizontalAlignment = -4108
$ws.Range("C3").VerticalAlignment = -4108

# Extend the sheet with a handful of blank rows below the data (formatted
# like the existing "mensagem" column cells) ready for new entries.
$ws.Range("D4:D7").HorizontalAlignment = -4108
$ws.Range("D4:D7").VerticalAlignment = -4108
$ws.Range("D4:D7").WrapText = $true

# Leave the selection where the user ended up after adding the new rows.
$ws.Range("D8").Select() | Out-Null
